{"js": "// The document contains a single 20-row x 5-column table. Only every\n// 4th row (0, 4, 8, 12, 16) actually holds data (the other rows are\n// spacer rows with empty paragraphs). Each populated row has 5 cells\n// with a \"A\u00f7B=C, D\" style division fact. This script replaces those\n// 25 values, by POSITION (table/row/column), with their new values -\n// this avoids any ambiguity from duplicate/overlapping text values\n// between the \"before\" and \"after\" sets.\nconst newValues = [\n  [\"29\u00f75=5, 4\", \"27\u00f76=4, 3\", \"54\u00f74=13, 2\", \"57\u00f78=7, 1\", \"84\u00f79=9, 3\"],\n  [\"19\u00f78=2, 3\", \"95\u00f79=10, 5\", \"58\u00f74=14, 2\", \"32\u00f77=4, 4\", \"89\u00f73=29, 2\"],\n  [\"52\u00f79=5, 7\", \"99\u00f79=11, 0\", \"62\u00f78=7, 6\", \"54\u00f77=7, 5\", \"25\u00f79=2, 7\"],\n  [\"54\u00f73=18, 0\", \"87\u00f73=29, 0\", \"36\u00f79=4, 0\", \"73\u00f76=12, 1\", \"55\u00f76=9, 1\"],\n  [\"77\u00f75=15, 2\", \"38\u00f76=6, 2\", \"95\u00f75=19, 0\", \"81\u00f79=9, 0\", \"50\u00f77=7, 1\"],\n];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const rowIdx = dataRowIndexes[i];\n  for (let col = 0; col < newValues[i].length; col++) {\n    const cell = table.getCell(rowIdx, col);\n    cell.value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table. Only every\n# 4th row (1, 5, 9, 13, 17 in 1-based Word indexing) actually holds\n# data (the other rows are spacer rows with empty paragraphs). Each\n# populated row has 5 cells with a \"A\u00f7B=C, D\" style division fact.\n# This script replaces those 25 values, by POSITION (table/row/col),\n# with their new values - this avoids any ambiguity from\n# duplicate/overlapping text values between the \"before\" and \"after\"\n# sets (one of the new values, 84\u00f79=9, 3, happens to equal one of the\n# *other* cells' original value, so a global text Find/Replace would\n# be unsafe here).\n\n$d = $word.ActiveDocument\n\n$table = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n    @(\"29\u00f75=5, 4\", \"27\u00f76=4, 3\", \"54\u00f74=13, 2\", \"57\u00f78=7, 1\", \"84\u00f79=9, 3\"),\n    @(\"19\u00f78=2, 3\", \"95\u00f79=10, 5\", \"58\u00f74=14, 2\", \"32\u00f77=4, 4\", \"89\u00f73=29, 2\"),\n    @(\"52\u00f79=5, 7\", \"99\u00f79=11, 0\", \"62\u00f78=7, 6\", \"54\u00f77=7, 5\", \"25\u00f79=2, 7\"),\n    @(\"54\u00f73=18, 0\", \"87\u00f73=29, 0\", \"36\u00f79=4, 0\", \"73\u00f76=12, 1\", \"55\u00f76=9, 1\"),\n    @(\"77\u00f75=15, 2\", \"38\u00f76=6, 2\", \"95\u00f75=19, 0\", \"81\u00f79=9, 0\", \"50\u00f77=7, 1\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $row = $dataRows[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $table.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
